$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "ISSFA - 0045"
$ws.Range("B8").Value = 44220
$ws.Range("B8").NumberFormat = $ws.Range("B7").NumberFormat
$ws.Range("C8").Value = "Iñaquito"
$ws.Range("D8").Value = "Quito"
$ws.Range("E8").Value = "Quito"
$ws.Range("F8").Value = "Quito"
$ws.Range("G8").Value = "Pichincha"
$ws.Range("H8").Value = "Casa"
$ws.Range("I8").Value = "Horizontal"
$ws.Range("J8").Value = 523
$ws.Range("K8").Value = 834
$ws.Range("L8").Value = 750.6
$ws.Range("M8").Value = 834
